$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 155, shifting existing rows (155..257) down to (156..258)
$ws.Rows.Item(155).Insert()

# Populate the newly inserted row 155 with the new record
$row = 155
$ws.Cells.Item($row, 1).Value  = 1
$ws.Cells.Item($row, 2).Value  = 'Agrícola del Norte S.A. de Arica'
$ws.Cells.Item($row, 3).Value  = 'Arica y Parinacota'
$ws.Cells.Item($row, 4).Value  = 45033
$ws.Cells.Item($row, 5).Value  = 15
$ws.Cells.Item($row, 6).Value  = 'Fruta'
$ws.Cells.Item($row, 7).Value  = 100108
$ws.Cells.Item($row, 8).Value  = 'Tropicales y subtropicales'
$ws.Cells.Item($row, 9).Value  = 100108002
$ws.Cells.Item($row, 10).Value = 'Mango'
$ws.Cells.Item($row, 11).Value = 'Sin especificar'
$ws.Cells.Item($row, 12).Value = 'Primera'
$ws.Cells.Item($row, 13).Value = 850
$ws.Cells.Item($row, 14).Value = 5000
$ws.Cells.Item($row, 15).Value = 5200
$ws.Cells.Item($row, 16).Value = 5071
$ws.Cells.Item($row, 17).Value = '$/bandeja 4 kilos'
$ws.Cells.Item($row, 18).Value = 'Perú'
$ws.Cells.Item($row, 19).Value = 1268
$ws.Cells.Item($row, 20).Value = 4
